$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.108.04"
$ws.Range("E2").Value = "  -1.27%  "

$ws.Range("D3").Value = "2.371.71"
$ws.Range("E3").Value = "  -4.69%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.79"
$ws.Range("E5").Value = "  -2.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.97"
$ws.Range("E6").Value = "  -5.10%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.526"
$ws.Range("E8").Value = "  -12.94%  "

$ws.Range("D9").Value = "2.369.84"
$ws.Range("E9").Value = "  -4.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.105"
$ws.Range("E10").Value = "  -3.77%  "

$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.28"
$ws.Range("E12").Value = "  -4.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.345"
$ws.Range("E13").Value = "  -4.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.16"
$ws.Range("E14").Value = "  -5.81%  "

$ws.Range("D15").Value = "2.804.97"
$ws.Range("E15").Value = "  -4.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000164"
$ws.Range("E16").Value = "  -3.64%  "

$ws.Range("D17").Value = "61.137.79"
$ws.Range("E17").Value = "  -1.05%  "

$ws.Range("D18").Value = "2.382.00"
$ws.Range("E18").Value = "  -3.98%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.68"
$ws.Range("E19").Value = "  -5.49%  "

$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.10"
$ws.Range("E20").Value = "  -3.75%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "317.27"
$ws.Range("E21").Value = "  -2.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.64"
$ws.Range("E22").Value = "  -8.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.15%  "

$ws.Range("E24").Value = "  -1.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.64"
$ws.Range("E25").Value = "  -1.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.15"
$ws.Range("E26").Value = "  +3.66%  "

$ws.Range("E27").Value = "  +0.29%  "

$ws.Range("D28").Value = "2.494.92"
$ws.Range("E28").Value = "  -4.28%  "

$ws.Range("D29").Value = "0.0₃0919"
$ws.Range("E29").Value = "  -10.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "523.43"
$ws.Range("E30").Value = "  -9.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.03"
$ws.Range("E31").Value = "  -4.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.41"
$ws.Range("E32").Value = "  -7.25%  "

$ws.Range("E33").Value = "  -3.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  -5.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("E35").Value = "  -3.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.63"
$ws.Range("E37").Value = "  -7.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.40"
$ws.Range("E38").Value = "  -10.44%  "

$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.374"
$ws.Range("E39").Value = "  -3.35%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.84"
$ws.Range("E40").Value = "  +2.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.11"
$ws.Range("E41").Value = "  -3.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.12"
$ws.Range("E42").Value = "  -3.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.41"
$ws.Range("E44").Value = "  -0.69%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "139.91"
$ws.Range("E45").Value = "  -6.73%  "

$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.12"
$ws.Range("E46").Value = "  -14.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.59"
$ws.Range("E47").Value = "  -2.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.89"
$ws.Range("E48").Value = "  -10.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0515"
$ws.Range("E49").Value = "  -5.89%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.572"
$ws.Range("E50").Value = "  -4.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0901"
$ws.Range("E51").Value = "  -5.12%  "
